$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row contents are permuted among rows 43/44, 46/47, 54/55, 56/57/58, 60/61
# (row numbers keep their position; the underlying record data moves between them).

# Row 43 <- data previously in row 44
$ws.Cells.Item(43, 1).Value = 130795187
$ws.Cells.Item(43, 2).Value = 79243
$ws.Cells.Item(43, 3).Value = $null
$ws.Cells.Item(43, 4).Value = "NT"
$ws.Cells.Item(43, 5).Value = 6425
$ws.Cells.Item(43, 6).Value = "Garnlav"
$ws.Cells.Item(43, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(43, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(43, 9).Value = "'30"
$ws.Cells.Item(43, 10).Value = "bålar"
$ws.Cells.Item(43, 11).Value = $null
$ws.Cells.Item(43, 12).Value = $null
$ws.Cells.Item(43, 13).Value = $null
$ws.Cells.Item(43, 14).Value = $null
$ws.Cells.Item(43, 15).Value = $null
$ws.Cells.Item(43, 16).Value = "Litnersmyren, Litnersmyren, Jmt"
$ws.Cells.Item(43, 17).Value = 444732
$ws.Cells.Item(43, 18).Value = 7025497
$ws.Cells.Item(43, 19).Value = 6
$ws.Cells.Item(43, 20).Value = "Jämtland"
$ws.Cells.Item(43, 21).Value = "Krokom"
$ws.Cells.Item(43, 22).Value = "Jämtland"
$ws.Cells.Item(43, 23).Value = "Alsen"
$ws.Cells.Item(43, 24).Value = $null
$ws.Cells.Item(43, 25).Value = "'2026-01-20"
$ws.Cells.Item(43, 26).Value = "12:21"
$ws.Cells.Item(43, 27).Value = "'2026-01-20"
$ws.Cells.Item(43, 28).Value = "12:21"
$ws.Cells.Item(43, 29).Value = "På gammal tall i gammal granskog"
$ws.Cells.Item(43, 30).Value = $false
$ws.Cells.Item(43, 31).Value = $false
$ws.Cells.Item(43, 32).Value = $null
$ws.Cells.Item(43, 33).Value = $false
$ws.Cells.Item(43, 34).Value = $null
$ws.Cells.Item(43, 35).Value = $null
$ws.Cells.Item(43, 36).Value = "tall"
$ws.Cells.Item(43, 37).Value = "Pinus sylvestris"
$ws.Cells.Item(43, 38).Value = $null
$ws.Cells.Item(43, 39).Value = $null
$ws.Cells.Item(43, 40).Value = $null
$ws.Cells.Item(43, 41).Value = "Pinus sylvestris"
$ws.Cells.Item(43, 42).Value = $null
$ws.Cells.Item(43, 43).Value = $null
$ws.Cells.Item(43, 44).Value = $null
$ws.Cells.Item(43, 45).Value = $null
$ws.Cells.Item(43, 46).Value = $null
$ws.Cells.Item(43, 47).Value = $null
$ws.Cells.Item(43, 48).Value = $null
$ws.Cells.Item(43, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(43, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(43, 51).Value = $null

# Row 44 <- data previously in row 43
$ws.Cells.Item(44, 1).Value = 130800374
$ws.Cells.Item(44, 2).Value = 83206
$ws.Cells.Item(44, 3).Value = $null
$ws.Cells.Item(44, 4).Value = "LC"
$ws.Cells.Item(44, 5).Value = 6439
$ws.Cells.Item(44, 6).Value = "Gulnål"
$ws.Cells.Item(44, 7).Value = "Chaenotheca brachypoda"
$ws.Cells.Item(44, 8).Value = "(Ach.) Tibell"
$ws.Cells.Item(44, 9).Value = $null
$ws.Cells.Item(44, 10).Value = $null
$ws.Cells.Item(44, 11).Value = $null
$ws.Cells.Item(44, 12).Value = $null
$ws.Cells.Item(44, 13).Value = $null
$ws.Cells.Item(44, 14).Value = $null
$ws.Cells.Item(44, 15).Value = $null
$ws.Cells.Item(44, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(44, 17).Value = 444662
$ws.Cells.Item(44, 18).Value = 7025556
$ws.Cells.Item(44, 19).Value = 10
$ws.Cells.Item(44, 20).Value = "Jämtland"
$ws.Cells.Item(44, 21).Value = "Krokom"
$ws.Cells.Item(44, 22).Value = "Jämtland"
$ws.Cells.Item(44, 23).Value = "Alsen"
$ws.Cells.Item(44, 24).Value = $null
$ws.Cells.Item(44, 25).Value = "'2026-01-20"
$ws.Cells.Item(44, 26).Value = "13:09"
$ws.Cells.Item(44, 27).Value = "'2026-01-20"
$ws.Cells.Item(44, 28).Value = "13:09"
$ws.Cells.Item(44, 29).Value = "Vid basen av grov björkhögstubbe (30 cm dbh) i gammal granskog med inslag av tallöverståndare"
$ws.Cells.Item(44, 30).Value = $false
$ws.Cells.Item(44, 31).Value = $false
$ws.Cells.Item(44, 32).Value = $null
$ws.Cells.Item(44, 33).Value = $false
$ws.Cells.Item(44, 34).Value = $null
$ws.Cells.Item(44, 35).Value = $null
$ws.Cells.Item(44, 36).Value = "glasbjörk"
$ws.Cells.Item(44, 37).Value = "Betula pubescens"
$ws.Cells.Item(44, 38).Value = $null
$ws.Cells.Item(44, 39).Value = $null
$ws.Cells.Item(44, 40).Value = $null
$ws.Cells.Item(44, 41).Value = "Betula pubescens"
$ws.Cells.Item(44, 42).Value = $null
$ws.Cells.Item(44, 43).Value = $null
$ws.Cells.Item(44, 44).Value = $null
$ws.Cells.Item(44, 45).Value = $null
$ws.Cells.Item(44, 46).Value = $null
$ws.Cells.Item(44, 47).Value = $null
$ws.Cells.Item(44, 48).Value = $null
$ws.Cells.Item(44, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(44, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(44, 51).Value = $null

# Row 46 <- data previously in row 47
$ws.Cells.Item(46, 1).Value = 130800372
$ws.Cells.Item(46, 2).Value = 79243
$ws.Cells.Item(46, 3).Value = $null
$ws.Cells.Item(46, 4).Value = "NT"
$ws.Cells.Item(46, 5).Value = 6425
$ws.Cells.Item(46, 6).Value = "Garnlav"
$ws.Cells.Item(46, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(46, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(46, 9).Value = $null
$ws.Cells.Item(46, 10).Value = $null
$ws.Cells.Item(46, 11).Value = $null
$ws.Cells.Item(46, 12).Value = $null
$ws.Cells.Item(46, 13).Value = $null
$ws.Cells.Item(46, 14).Value = $null
$ws.Cells.Item(46, 15).Value = $null
$ws.Cells.Item(46, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(46, 17).Value = 444688
$ws.Cells.Item(46, 18).Value = 7025582
$ws.Cells.Item(46, 19).Value = 10
$ws.Cells.Item(46, 20).Value = "Jämtland"
$ws.Cells.Item(46, 21).Value = "Krokom"
$ws.Cells.Item(46, 22).Value = "Jämtland"
$ws.Cells.Item(46, 23).Value = "Alsen"
$ws.Cells.Item(46, 24).Value = $null
$ws.Cells.Item(46, 25).Value = "'2026-01-20"
$ws.Cells.Item(46, 26).Value = "13:20"
$ws.Cells.Item(46, 27).Value = "'2026-01-20"
$ws.Cells.Item(46, 28).Value = "13:20"
$ws.Cells.Item(46, 29).Value = "På gammal gran i gammal granskog"
$ws.Cells.Item(46, 30).Value = $false
$ws.Cells.Item(46, 31).Value = $false
$ws.Cells.Item(46, 32).Value = $null
$ws.Cells.Item(46, 33).Value = $false
$ws.Cells.Item(46, 34).Value = $null
$ws.Cells.Item(46, 35).Value = $null
$ws.Cells.Item(46, 36).Value = "gran"
$ws.Cells.Item(46, 37).Value = "Picea abies"
$ws.Cells.Item(46, 38).Value = $null
$ws.Cells.Item(46, 39).Value = $null
$ws.Cells.Item(46, 40).Value = $null
$ws.Cells.Item(46, 41).Value = "Picea abies"
$ws.Cells.Item(46, 42).Value = $null
$ws.Cells.Item(46, 43).Value = $null
$ws.Cells.Item(46, 44).Value = $null
$ws.Cells.Item(46, 45).Value = $null
$ws.Cells.Item(46, 46).Value = $null
$ws.Cells.Item(46, 47).Value = $null
$ws.Cells.Item(46, 48).Value = $null
$ws.Cells.Item(46, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(46, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(46, 51).Value = $null

# Row 47 <- data previously in row 46
$ws.Cells.Item(47, 1).Value = 130795383
$ws.Cells.Item(47, 2).Value = 80348
$ws.Cells.Item(47, 3).Value = $null
$ws.Cells.Item(47, 4).Value = "NT"
$ws.Cells.Item(47, 5).Value = 6458
$ws.Cells.Item(47, 6).Value = "Lunglav"
$ws.Cells.Item(47, 7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(47, 8).Value = "(L.) Hoffm."
$ws.Cells.Item(47, 9).Value = $null
$ws.Cells.Item(47, 10).Value = $null
$ws.Cells.Item(47, 11).Value = $null
$ws.Cells.Item(47, 12).Value = $null
$ws.Cells.Item(47, 13).Value = $null
$ws.Cells.Item(47, 14).Value = $null
$ws.Cells.Item(47, 15).Value = $null
$ws.Cells.Item(47, 16).Value = "Litnersmyren, Litnersmyren, Jmt"
$ws.Cells.Item(47, 17).Value = 444726
$ws.Cells.Item(47, 18).Value = 7025453
$ws.Cells.Item(47, 19).Value = 10
$ws.Cells.Item(47, 20).Value = "Jämtland"
$ws.Cells.Item(47, 21).Value = "Krokom"
$ws.Cells.Item(47, 22).Value = "Jämtland"
$ws.Cells.Item(47, 23).Value = "Alsen"
$ws.Cells.Item(47, 24).Value = $null
$ws.Cells.Item(47, 25).Value = "'2026-01-20"
$ws.Cells.Item(47, 26).Value = "12:32"
$ws.Cells.Item(47, 27).Value = "'2026-01-20"
$ws.Cells.Item(47, 28).Value = "12:32"
$ws.Cells.Item(47, 29).Value = "På gammal klen sälg i gammal barrblandskog"
$ws.Cells.Item(47, 30).Value = $false
$ws.Cells.Item(47, 31).Value = $false
$ws.Cells.Item(47, 32).Value = $null
$ws.Cells.Item(47, 33).Value = $false
$ws.Cells.Item(47, 34).Value = $null
$ws.Cells.Item(47, 35).Value = $null
$ws.Cells.Item(47, 36).Value = "sälg"
$ws.Cells.Item(47, 37).Value = "Salix caprea"
$ws.Cells.Item(47, 38).Value = $null
$ws.Cells.Item(47, 39).Value = $null
$ws.Cells.Item(47, 40).Value = $null
$ws.Cells.Item(47, 41).Value = "Salix caprea"
$ws.Cells.Item(47, 42).Value = $null
$ws.Cells.Item(47, 43).Value = $null
$ws.Cells.Item(47, 44).Value = $null
$ws.Cells.Item(47, 45).Value = $null
$ws.Cells.Item(47, 46).Value = $null
$ws.Cells.Item(47, 47).Value = $null
$ws.Cells.Item(47, 48).Value = $null
$ws.Cells.Item(47, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(47, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(47, 51).Value = $null

# Row 54 <- data previously in row 55
$ws.Cells.Item(54, 1).Value = 130800375
$ws.Cells.Item(54, 2).Value = 83214
$ws.Cells.Item(54, 3).Value = $null
$ws.Cells.Item(54, 4).Value = "VU"
$ws.Cells.Item(54, 5).Value = 492
$ws.Cells.Item(54, 6).Value = "Smalskaftslav"
$ws.Cells.Item(54, 7).Value = "Chaenotheca gracilenta"
$ws.Cells.Item(54, 8).Value = "(Ach.) J.Mattsson & Middelb."
$ws.Cells.Item(54, 9).Value = $null
$ws.Cells.Item(54, 10).Value = $null
$ws.Cells.Item(54, 11).Value = $null
$ws.Cells.Item(54, 12).Value = $null
$ws.Cells.Item(54, 13).Value = $null
$ws.Cells.Item(54, 14).Value = $null
$ws.Cells.Item(54, 15).Value = $null
$ws.Cells.Item(54, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(54, 17).Value = 444662
$ws.Cells.Item(54, 18).Value = 7025556
$ws.Cells.Item(54, 19).Value = 10
$ws.Cells.Item(54, 20).Value = "Jämtland"
$ws.Cells.Item(54, 21).Value = "Krokom"
$ws.Cells.Item(54, 22).Value = "Jämtland"
$ws.Cells.Item(54, 23).Value = "Alsen"
$ws.Cells.Item(54, 24).Value = $null
$ws.Cells.Item(54, 25).Value = "'2026-01-20"
$ws.Cells.Item(54, 26).Value = "13:09"
$ws.Cells.Item(54, 27).Value = "'2026-01-20"
$ws.Cells.Item(54, 28).Value = "13:09"
$ws.Cells.Item(54, 29).Value = "Vid basen av grov björkhögstubbe (30 cm dbh) i gammal granskog med inslag av tallöverståndare"
$ws.Cells.Item(54, 30).Value = $false
$ws.Cells.Item(54, 31).Value = $false
$ws.Cells.Item(54, 32).Value = $null
$ws.Cells.Item(54, 33).Value = $false
$ws.Cells.Item(54, 34).Value = $null
$ws.Cells.Item(54, 35).Value = $null
$ws.Cells.Item(54, 36).Value = "glasbjörk"
$ws.Cells.Item(54, 37).Value = "Betula pubescens"
$ws.Cells.Item(54, 38).Value = $null
$ws.Cells.Item(54, 39).Value = $null
$ws.Cells.Item(54, 40).Value = $null
$ws.Cells.Item(54, 41).Value = "Betula pubescens"
$ws.Cells.Item(54, 42).Value = $null
$ws.Cells.Item(54, 43).Value = $null
$ws.Cells.Item(54, 44).Value = $null
$ws.Cells.Item(54, 45).Value = $null
$ws.Cells.Item(54, 46).Value = $null
$ws.Cells.Item(54, 47).Value = $null
$ws.Cells.Item(54, 48).Value = $null
$ws.Cells.Item(54, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(54, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(54, 51).Value = $null

# Row 55 <- data previously in row 54
$ws.Cells.Item(55, 1).Value = 130800352
$ws.Cells.Item(55, 2).Value = 79243
$ws.Cells.Item(55, 3).Value = $null
$ws.Cells.Item(55, 4).Value = "NT"
$ws.Cells.Item(55, 5).Value = 6425
$ws.Cells.Item(55, 6).Value = "Garnlav"
$ws.Cells.Item(55, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(55, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(55, 9).Value = $null
$ws.Cells.Item(55, 10).Value = $null
$ws.Cells.Item(55, 11).Value = $null
$ws.Cells.Item(55, 12).Value = $null
$ws.Cells.Item(55, 13).Value = $null
$ws.Cells.Item(55, 14).Value = $null
$ws.Cells.Item(55, 15).Value = $null
$ws.Cells.Item(55, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(55, 17).Value = 444700
$ws.Cells.Item(55, 18).Value = 7025517
$ws.Cells.Item(55, 19).Value = 10
$ws.Cells.Item(55, 20).Value = "Jämtland"
$ws.Cells.Item(55, 21).Value = "Krokom"
$ws.Cells.Item(55, 22).Value = "Jämtland"
$ws.Cells.Item(55, 23).Value = "Alsen"
$ws.Cells.Item(55, 24).Value = $null
$ws.Cells.Item(55, 25).Value = "'2026-01-20"
$ws.Cells.Item(55, 26).Value = "14:43"
$ws.Cells.Item(55, 27).Value = "'2026-01-20"
$ws.Cells.Item(55, 28).Value = "14:43"
$ws.Cells.Item(55, 29).Value = "Rikligt på gammal levande gran i gammal granskog"
$ws.Cells.Item(55, 30).Value = $false
$ws.Cells.Item(55, 31).Value = $false
$ws.Cells.Item(55, 32).Value = $null
$ws.Cells.Item(55, 33).Value = $false
$ws.Cells.Item(55, 34).Value = $null
$ws.Cells.Item(55, 35).Value = $null
$ws.Cells.Item(55, 36).Value = "gran"
$ws.Cells.Item(55, 37).Value = "Picea abies"
$ws.Cells.Item(55, 38).Value = $null
$ws.Cells.Item(55, 39).Value = $null
$ws.Cells.Item(55, 40).Value = $null
$ws.Cells.Item(55, 41).Value = "Picea abies"
$ws.Cells.Item(55, 42).Value = $null
$ws.Cells.Item(55, 43).Value = $null
$ws.Cells.Item(55, 44).Value = $null
$ws.Cells.Item(55, 45).Value = $null
$ws.Cells.Item(55, 46).Value = $null
$ws.Cells.Item(55, 47).Value = $null
$ws.Cells.Item(55, 48).Value = $null
$ws.Cells.Item(55, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(55, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(55, 51).Value = $null

# Row 56 <- data previously in row 57
$ws.Cells.Item(56, 1).Value = 130800361
$ws.Cells.Item(56, 2).Value = 75333
$ws.Cells.Item(56, 3).Value = $null
$ws.Cells.Item(56, 4).Value = "NT"
$ws.Cells.Item(56, 5).Value = 1460
$ws.Cells.Item(56, 6).Value = "Rosa skärelav"
$ws.Cells.Item(56, 7).Value = "Schismatomma pericleum"
$ws.Cells.Item(56, 8).Value = "(Ach.) Branth & Rostr."
$ws.Cells.Item(56, 9).Value = $null
$ws.Cells.Item(56, 10).Value = $null
$ws.Cells.Item(56, 11).Value = $null
$ws.Cells.Item(56, 12).Value = $null
$ws.Cells.Item(56, 13).Value = $null
$ws.Cells.Item(56, 14).Value = $null
$ws.Cells.Item(56, 15).Value = $null
$ws.Cells.Item(56, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(56, 17).Value = 444731
$ws.Cells.Item(56, 18).Value = 7025516
$ws.Cells.Item(56, 19).Value = 10
$ws.Cells.Item(56, 20).Value = "Jämtland"
$ws.Cells.Item(56, 21).Value = "Krokom"
$ws.Cells.Item(56, 22).Value = "Jämtland"
$ws.Cells.Item(56, 23).Value = "Alsen"
$ws.Cells.Item(56, 24).Value = $null
$ws.Cells.Item(56, 25).Value = "'2026-01-20"
$ws.Cells.Item(56, 26).Value = "13:54"
$ws.Cells.Item(56, 27).Value = "'2026-01-20"
$ws.Cells.Item(56, 28).Value = "13:54"
$ws.Cells.Item(56, 29).Value = "Vid basen av grov gammal levande sälg (50 cm dbh) i gammal granskog"
$ws.Cells.Item(56, 30).Value = $false
$ws.Cells.Item(56, 31).Value = $false
$ws.Cells.Item(56, 32).Value = $null
$ws.Cells.Item(56, 33).Value = $false
$ws.Cells.Item(56, 34).Value = $null
$ws.Cells.Item(56, 35).Value = $null
$ws.Cells.Item(56, 36).Value = "sälg"
$ws.Cells.Item(56, 37).Value = "Salix caprea"
$ws.Cells.Item(56, 38).Value = $null
$ws.Cells.Item(56, 39).Value = $null
$ws.Cells.Item(56, 40).Value = $null
$ws.Cells.Item(56, 41).Value = "Salix caprea"
$ws.Cells.Item(56, 42).Value = $null
$ws.Cells.Item(56, 43).Value = $null
$ws.Cells.Item(56, 44).Value = $null
$ws.Cells.Item(56, 45).Value = $null
$ws.Cells.Item(56, 46).Value = $null
$ws.Cells.Item(56, 47).Value = $null
$ws.Cells.Item(56, 48).Value = $null
$ws.Cells.Item(56, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(56, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(56, 51).Value = $null

# Row 57 <- data previously in row 58
$ws.Cells.Item(57, 1).Value = 130800364
$ws.Cells.Item(57, 2).Value = 80384
$ws.Cells.Item(57, 3).Value = $null
$ws.Cells.Item(57, 4).Value = "LC"
$ws.Cells.Item(57, 5).Value = 6464
$ws.Cells.Item(57, 6).Value = "Luddlav"
$ws.Cells.Item(57, 7).Value = "Nephroma resupinatum"
$ws.Cells.Item(57, 8).Value = "(L.) Ach."
$ws.Cells.Item(57, 9).Value = $null
$ws.Cells.Item(57, 10).Value = $null
$ws.Cells.Item(57, 11).Value = $null
$ws.Cells.Item(57, 12).Value = $null
$ws.Cells.Item(57, 13).Value = $null
$ws.Cells.Item(57, 14).Value = $null
$ws.Cells.Item(57, 15).Value = $null
$ws.Cells.Item(57, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(57, 17).Value = 444717
$ws.Cells.Item(57, 18).Value = 7025526
$ws.Cells.Item(57, 19).Value = 10
$ws.Cells.Item(57, 20).Value = "Jämtland"
$ws.Cells.Item(57, 21).Value = "Krokom"
$ws.Cells.Item(57, 22).Value = "Jämtland"
$ws.Cells.Item(57, 23).Value = "Alsen"
$ws.Cells.Item(57, 24).Value = $null
$ws.Cells.Item(57, 25).Value = "'2026-01-20"
$ws.Cells.Item(57, 26).Value = "13:48"
$ws.Cells.Item(57, 27).Value = "'2026-01-20"
$ws.Cells.Item(57, 28).Value = "13:48"
$ws.Cells.Item(57, 29).Value = "På bark på stam av levande lutande gammal sälg i gammal granskog"
$ws.Cells.Item(57, 30).Value = $false
$ws.Cells.Item(57, 31).Value = $false
$ws.Cells.Item(57, 32).Value = $null
$ws.Cells.Item(57, 33).Value = $false
$ws.Cells.Item(57, 34).Value = $null
$ws.Cells.Item(57, 35).Value = $null
$ws.Cells.Item(57, 36).Value = "sälg"
$ws.Cells.Item(57, 37).Value = "Salix caprea"
$ws.Cells.Item(57, 38).Value = $null
$ws.Cells.Item(57, 39).Value = $null
$ws.Cells.Item(57, 40).Value = $null
$ws.Cells.Item(57, 41).Value = "Salix caprea"
$ws.Cells.Item(57, 42).Value = $null
$ws.Cells.Item(57, 43).Value = $null
$ws.Cells.Item(57, 44).Value = $null
$ws.Cells.Item(57, 45).Value = $null
$ws.Cells.Item(57, 46).Value = $null
$ws.Cells.Item(57, 47).Value = $null
$ws.Cells.Item(57, 48).Value = $null
$ws.Cells.Item(57, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(57, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(57, 51).Value = $null

# Row 58 <- data previously in row 56
$ws.Cells.Item(58, 1).Value = 130800354
$ws.Cells.Item(58, 2).Value = 79243
$ws.Cells.Item(58, 3).Value = $null
$ws.Cells.Item(58, 4).Value = "NT"
$ws.Cells.Item(58, 5).Value = 6425
$ws.Cells.Item(58, 6).Value = "Garnlav"
$ws.Cells.Item(58, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(58, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(58, 9).Value = $null
$ws.Cells.Item(58, 10).Value = $null
$ws.Cells.Item(58, 11).Value = $null
$ws.Cells.Item(58, 12).Value = $null
$ws.Cells.Item(58, 13).Value = $null
$ws.Cells.Item(58, 14).Value = $null
$ws.Cells.Item(58, 15).Value = $null
$ws.Cells.Item(58, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(58, 17).Value = 444716
$ws.Cells.Item(58, 18).Value = 7025439
$ws.Cells.Item(58, 19).Value = 10
$ws.Cells.Item(58, 20).Value = "Jämtland"
$ws.Cells.Item(58, 21).Value = "Krokom"
$ws.Cells.Item(58, 22).Value = "Jämtland"
$ws.Cells.Item(58, 23).Value = "Alsen"
$ws.Cells.Item(58, 24).Value = $null
$ws.Cells.Item(58, 25).Value = "'2026-01-20"
$ws.Cells.Item(58, 26).Value = "14:16"
$ws.Cells.Item(58, 27).Value = "'2026-01-20"
$ws.Cells.Item(58, 28).Value = "14:16"
$ws.Cells.Item(58, 29).Value = "På gammal gran i gammal granskog"
$ws.Cells.Item(58, 30).Value = $false
$ws.Cells.Item(58, 31).Value = $false
$ws.Cells.Item(58, 32).Value = $null
$ws.Cells.Item(58, 33).Value = $false
$ws.Cells.Item(58, 34).Value = $null
$ws.Cells.Item(58, 35).Value = $null
$ws.Cells.Item(58, 36).Value = "gran"
$ws.Cells.Item(58, 37).Value = "Picea abies"
$ws.Cells.Item(58, 38).Value = $null
$ws.Cells.Item(58, 39).Value = $null
$ws.Cells.Item(58, 40).Value = $null
$ws.Cells.Item(58, 41).Value = "Picea abies"
$ws.Cells.Item(58, 42).Value = $null
$ws.Cells.Item(58, 43).Value = $null
$ws.Cells.Item(58, 44).Value = $null
$ws.Cells.Item(58, 45).Value = $null
$ws.Cells.Item(58, 46).Value = $null
$ws.Cells.Item(58, 47).Value = $null
$ws.Cells.Item(58, 48).Value = $null
$ws.Cells.Item(58, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(58, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(58, 51).Value = $null

# Row 60 <- data previously in row 61
$ws.Cells.Item(60, 1).Value = 130800371
$ws.Cells.Item(60, 2).Value = 78255
$ws.Cells.Item(60, 3).Value = $null
$ws.Cells.Item(60, 4).Value = "NT"
$ws.Cells.Item(60, 5).Value = 228579
$ws.Cells.Item(60, 6).Value = "Liten svartspik"
$ws.Cells.Item(60, 7).Value = "Chaenothecopsis nana"
$ws.Cells.Item(60, 8).Value = "Tibell"
$ws.Cells.Item(60, 9).Value = $null
$ws.Cells.Item(60, 10).Value = $null
$ws.Cells.Item(60, 11).Value = $null
$ws.Cells.Item(60, 12).Value = $null
$ws.Cells.Item(60, 13).Value = $null
$ws.Cells.Item(60, 14).Value = $null
$ws.Cells.Item(60, 15).Value = $null
$ws.Cells.Item(60, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(60, 17).Value = 444684
$ws.Cells.Item(60, 18).Value = 7025581
$ws.Cells.Item(60, 19).Value = 10
$ws.Cells.Item(60, 20).Value = "Jämtland"
$ws.Cells.Item(60, 21).Value = "Krokom"
$ws.Cells.Item(60, 22).Value = "Jämtland"
$ws.Cells.Item(60, 23).Value = "Alsen"
$ws.Cells.Item(60, 24).Value = $null
$ws.Cells.Item(60, 25).Value = "'2026-01-20"
$ws.Cells.Item(60, 26).Value = "13:22"
$ws.Cells.Item(60, 27).Value = "'2026-01-20"
$ws.Cells.Item(60, 28).Value = "13:22"
$ws.Cells.Item(60, 29).Value = "På bark på stam av levande gammal gran i gammal granskog"
$ws.Cells.Item(60, 30).Value = $false
$ws.Cells.Item(60, 31).Value = $false
$ws.Cells.Item(60, 32).Value = $null
$ws.Cells.Item(60, 33).Value = $false
$ws.Cells.Item(60, 34).Value = $null
$ws.Cells.Item(60, 35).Value = $null
$ws.Cells.Item(60, 36).Value = "gran"
$ws.Cells.Item(60, 37).Value = "Picea abies"
$ws.Cells.Item(60, 38).Value = $null
$ws.Cells.Item(60, 39).Value = $null
$ws.Cells.Item(60, 40).Value = $null
$ws.Cells.Item(60, 41).Value = "Picea abies"
$ws.Cells.Item(60, 42).Value = $null
$ws.Cells.Item(60, 43).Value = $null
$ws.Cells.Item(60, 44).Value = $null
$ws.Cells.Item(60, 45).Value = $null
$ws.Cells.Item(60, 46).Value = $null
$ws.Cells.Item(60, 47).Value = $null
$ws.Cells.Item(60, 48).Value = $null
$ws.Cells.Item(60, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(60, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(60, 51).Value = $null

# Row 61 <- data previously in row 60
$ws.Cells.Item(61, 1).Value = 130800356
$ws.Cells.Item(61, 2).Value = 80349
$ws.Cells.Item(61, 3).Value = $null
$ws.Cells.Item(61, 4).Value = "NT"
$ws.Cells.Item(61, 5).Value = 2081
$ws.Cells.Item(61, 6).Value = "Skrovellav"
$ws.Cells.Item(61, 7).Value = "Lobaria scrobiculata"
$ws.Cells.Item(61, 8).Value = "(Scop.) DC."
$ws.Cells.Item(61, 9).Value = $null
$ws.Cells.Item(61, 10).Value = $null
$ws.Cells.Item(61, 11).Value = $null
$ws.Cells.Item(61, 12).Value = $null
$ws.Cells.Item(61, 13).Value = $null
$ws.Cells.Item(61, 14).Value = $null
$ws.Cells.Item(61, 15).Value = $null
$ws.Cells.Item(61, 16).Value = "Litnersmyran, Jmt"
$ws.Cells.Item(61, 17).Value = 444731
$ws.Cells.Item(61, 18).Value = 7025516
$ws.Cells.Item(61, 19).Value = 10
$ws.Cells.Item(61, 20).Value = "Jämtland"
$ws.Cells.Item(61, 21).Value = "Krokom"
$ws.Cells.Item(61, 22).Value = "Jämtland"
$ws.Cells.Item(61, 23).Value = "Alsen"
$ws.Cells.Item(61, 24).Value = $null
$ws.Cells.Item(61, 25).Value = "'2026-01-20"
$ws.Cells.Item(61, 26).Value = "13:54"
$ws.Cells.Item(61, 27).Value = "'2026-01-20"
$ws.Cells.Item(61, 28).Value = "13:54"
$ws.Cells.Item(61, 29).Value = "På bark av grov gammal levande sälg (50 cm dbh) i gammal granskog"
$ws.Cells.Item(61, 30).Value = $false
$ws.Cells.Item(61, 31).Value = $false
$ws.Cells.Item(61, 32).Value = $null
$ws.Cells.Item(61, 33).Value = $false
$ws.Cells.Item(61, 34).Value = $null
$ws.Cells.Item(61, 35).Value = $null
$ws.Cells.Item(61, 36).Value = "sälg"
$ws.Cells.Item(61, 37).Value = "Salix caprea"
$ws.Cells.Item(61, 38).Value = $null
$ws.Cells.Item(61, 39).Value = $null
$ws.Cells.Item(61, 40).Value = $null
$ws.Cells.Item(61, 41).Value = "Salix caprea"
$ws.Cells.Item(61, 42).Value = $null
$ws.Cells.Item(61, 43).Value = $null
$ws.Cells.Item(61, 44).Value = $null
$ws.Cells.Item(61, 45).Value = $null
$ws.Cells.Item(61, 46).Value = $null
$ws.Cells.Item(61, 47).Value = $null
$ws.Cells.Item(61, 48).Value = $null
$ws.Cells.Item(61, 49).Value = "Fredrik Jonsson"
$ws.Cells.Item(61, 50).Value = "Fredrik Jonsson"
$ws.Cells.Item(61, 51).Value = $null
